$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 306 (shifts existing rows 306-379 down to 307-380)
$ws.Rows.Item(306).Insert(-4121)

# Populate the newly inserted row 306 with the new record's data
$ws.Cells.Item(306, 1).Value = 4
$ws.Cells.Item(306, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(306, 3).Value = "Los Lagos"
$ws.Cells.Item(306, 4).Value = 44782
$ws.Cells.Item(306, 5).Value = 10
$ws.Cells.Item(306, 6).Value = 100112023
$ws.Cells.Item(306, 7).Value = "Brócoli"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 1400
$ws.Cells.Item(306, 11).Value = 1500
$ws.Cells.Item(306, 12).Value = 1500
$ws.Cells.Item(306, 13).Value = 1500
$ws.Cells.Item(306, 14).Value = "$/unidad"
$ws.Cells.Item(306, 15).Value = "Región Metropolitana"
$ws.Cells.Item(306, 16).Value = 1500
$ws.Cells.Item(306, 17).Value = 1
$ws.Cells.Item(306, 18).Value = "Hortaliza"
